$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "jonas"
$ws.Range("C4").Value = "karlsson"
$ws.Range("D4").Value = "mail@mail.com"
$ws.Range("E4").Value = 187
$ws.Range("F4").Value = 205
$ws.Range("G4").Value = 44
